$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Errors")

$ws.Range("B7").Value = "Enemy clipping and the height"
$ws.Range("B8").Value = "Also scaling moves ui"

$ws.Range("G4").Value = "Add by tomo"
$ws.Range("G5").Value = "Audio"
$ws.Range("G6").Value = "UI"
$ws.Range("G7").Value = "Enemy working"
$ws.Range("G8").Value = "Game end"
$ws.Range("G9").Value = "Boss working"

$ws.Activate()
$ws.Range("G10").Select()
